$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 11 - shifts existing rows 11:112 down to 12:113
$ws.Rows("11:11").Insert()

# Populate the new row 11 with its data
$ws.Range("A11").Value = 11
$ws.Range("B11").Value = "Vega Monumental Concepción"
$ws.Range("C11").Value = "Bíobío"
$ws.Range("D11").Value = 44649
$ws.Range("E11").Value = 8
$ws.Range("F11").Value = 100112032
$ws.Range("G11").Value = "Zapallo italiano"
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 220
$ws.Range("K11").Value = 12000
$ws.Range("L11").Value = 13000
$ws.Range("M11").Value = 12455
$ws.Range("N11").Value = "$/caja 60 unidades"
$ws.Range("O11").Value = "Región de Arica y Parinacota"
$ws.Range("P11").Value = 208
$ws.Range("Q11").Value = 60
$ws.Range("R11").Value = "Hortaliza"

# Apply the same date number format/style as the other D-column cells
$ws.Range("D11").NumberFormat = $ws.Range("D12").NumberFormat
